# Presence of properties in articles.xlsx - remove the per-section breakdown
# columns (D:P) from Sheet2, leaving just PII / Table / Text, and move the
# active selection to G12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Remove all the per-section columns (D through P) for the header row and
# every data row (2-51). Using Clear() (rather than ClearContents()) drops
# the cells entirely -- including their style -- matching the target file,
# which has no <c> elements at all for these columns anymore. Column Q is
# intentionally excluded so the Q12 placeholder cell is left untouched.
$ws.Range("D1:P51").Clear()

# The header row no longer needs the explicit 45pt row height now that the
# wrapped section-name headers are gone; let Excel recompute the default.
$ws.Rows.Item(1).AutoFit()

# Update the saved cursor/selection position.
$ws.Range("G12").Select()

Write-Host "Removed section-breakdown columns from Sheet2 and updated selection."
